# Location filter, Export deal, Negotiations and Chat Connection test cases updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$data = @(
    @("LocFilter_TC001", "John Tucker", "ONE"),
    @("LocFilter_TC002", "NA", "ALL"),
    @("LocFilter_TC003", "Stan Koster Andersons", "ONE"),
    @("NegotiateDeal_TC001", "John Tucker", "ONE"),
    @("NegotiateDeal_TC002", "NA", "ALL"),
    @("NegotiateDeal_TC003", "Stan Koster Andersons", "ONE"),
    @("NegotiateDeal_TC004", "John Tucker", "ONE")
)

$row = 31
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = "Deal shared successfully"
    $row = $row + 1
}

$ws.Range("A31:A37").VerticalAlignment = -4108

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A31:D37").Select()
